$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(3, 6).Value = 54797
$ws.Cells.Item(4, 6).Value = 1345
$ws.Cells.Item(5, 6).Value = 387
$ws.Cells.Item(6, 6).Value = 328
$ws.Cells.Item(7, 6).Value = 884
$ws.Cells.Item(8, 6).Value = 765
$ws.Cells.Item(9, 6).Value = 408
$ws.Cells.Item(10, 6).Value = 3071
$ws.Cells.Item(11, 6).Value = 913
$ws.Cells.Item(12, 6).Value = 5237
$ws.Cells.Item(13, 6).Value = 1286
$ws.Cells.Item(13, 7).Value = 125
$ws.Cells.Item(14, 6).Value = 1051
$ws.Cells.Item(16, 6).Value = 847
$ws.Cells.Item(18, 6).Value = 411
$ws.Cells.Item(19, 6).Value = 1294
$ws.Cells.Item(20, 6).Value = 104
$ws.Cells.Item(22, 6).Value = 186
$ws.Cells.Item(23, 6).Value = 370
$ws.Cells.Item(24, 6).Value = 32
$ws.Cells.Item(25, 6).Value = 39
$ws.Cells.Item(27, 6).Value = 68
$ws.Cells.Item(28, 6).Value = 62
$ws.Cells.Item(29, 6).Value = 5111
$ws.Cells.Item(30, 6).Value = 37
$ws.Cells.Item(31, 6).Value = 5015
$ws.Cells.Item(32, 6).Value = 9022
$ws.Cells.Item(33, 6).Value = 114
$ws.Cells.Item(34, 6).Value = 153
$ws.Cells.Item(36, 6).Value = 227
$ws.Cells.Item(37, 6).Value = 433
$ws.Cells.Item(38, 6).Value = 117
$ws.Cells.Item(39, 6).Value = 87
$ws.Cells.Item(40, 6).Value = 4217
$ws.Cells.Item(41, 6).Value = 247

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(4, 6).Value = 95
$ws.Cells.Item(12, 6).Value = 1134

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 785
$ws.Cells.Item(3, 6).Value = 574
$ws.Cells.Item(4, 6).Value = 140
$ws.Cells.Item(5, 6).Value = 43

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 785
$ws.Cells.Item(3, 6).Value = 1345
$ws.Cells.Item(4, 6).Value = 387
$ws.Cells.Item(5, 6).Value = 328
$ws.Cells.Item(6, 6).Value = 884
$ws.Cells.Item(7, 6).Value = 765
$ws.Cells.Item(8, 6).Value = 408
$ws.Cells.Item(9, 6).Value = 913
$ws.Cells.Item(10, 6).Value = 95
$ws.Cells.Item(11, 6).Value = 1286
$ws.Cells.Item(11, 7).Value = 125
$ws.Cells.Item(12, 6).Value = 43
$ws.Cells.Item(14, 6).Value = 1051
$ws.Cells.Item(16, 6).Value = 847
$ws.Cells.Item(17, 6).Value = 411
$ws.Cells.Item(19, 6).Value = 1294
$ws.Cells.Item(21, 6).Value = 104
$ws.Cells.Item(22, 6).Value = 186
$ws.Cells.Item(24, 6).Value = 370
$ws.Cells.Item(25, 6).Value = 32
$ws.Cells.Item(26, 6).Value = 39
$ws.Cells.Item(27, 6).Value = 62
$ws.Cells.Item(28, 6).Value = 5111
$ws.Cells.Item(29, 6).Value = 37
$ws.Cells.Item(30, 6).Value = 9022
$ws.Cells.Item(32, 6).Value = 114
$ws.Cells.Item(33, 6).Value = 153
$ws.Cells.Item(34, 6).Value = 137
$ws.Cells.Item(35, 6).Value = 227
$ws.Cells.Item(36, 6).Value = 433
$ws.Cells.Item(39, 6).Value = 117
$ws.Cells.Item(40, 6).Value = 87
$ws.Cells.Item(41, 6).Value = 4217
$ws.Cells.Item(48, 6).Value = 247
